$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the header/label column with a new "coefficients table" wording:
# A1 becomes the column-header label "var"; B1/C1/D1 keep their HBW/HBO/NHB labels.
$ws.Range("A1").Value = "var"
$ws.Range("B1").Value = "HBW"
$ws.Range("C1").Value = "HBO"
$ws.Range("D1").Value = "NHB"

# Expand the abbreviated row labels in column A into full descriptive names.
$ws.Range("A2").Value = "Households"
$ws.Range("A3").Value = "Office Employment"
$ws.Range("A4").Value = "Other Employment"
$ws.Range("A5").Value = "Retail Employment"

# Widen column A so the longer labels are fully visible.
$ws.Columns.Item(1).ColumnWidth = 21.333333333333332

# Leave the cursor on the row below the table, as in the saved file.
$ws.Range("A6").Select() | Out-Null
